$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for a new row at 207 by shifting rows 207-214 down to 208-215 ---
# Copy formatting first (bottom-up) so every destination row gets the exact
# formatting of the row that is moving into it, without Excel synthesizing
# brand-new style combinations.
for ($r = 214; $r -ge 207; $r--) {
    $dst = $r + 1
    $ws.Range("A$r`:F$r").Copy() | Out-Null
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Now move the values themselves (bottom-up so we never clobber a row before
# reading it).
for ($r = 214; $r -ge 207; $r--) {
    $dst = $r + 1
    for ($c = 1; $c -le 6; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dst, $c)
        $dstCell.Value2 = $srcCell.Value2
    }
}

Write-Host "row shift complete"

# --- Step 2: populate the newly freed row 207 with the humidity indicator ---
$ws.Range("A207").Value2 = 205
$ws.Range("B207").Value2 = "main_indicator_humidity"
$ws.Range("C207").Value2 = "Humidity"
$ws.Range("D207").Value2 = "濕度"
$ws.Range("E207").Value2 = "Feuchtigkeit"
$ws.Range("F207").Value2 = "湿度"

# --- Step 3: renumber column A for rows 208-215 (shifted by one) ---
for ($r = 208; $r -le 215; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

Write-Host "row207 + renumber complete"

# --- Step 4: append the new sequence-editor / teaching-guide rows at the bottom ---
$ws.Range("A215:F215").Copy() | Out-Null
$ws.Range("A216:F222").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A216:F222").Value2 = $ws.Range("A209:F215").Value2
for ($r = 216; $r -le 222; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

$ws.Range("B216").Value2 = "seq_step"
$ws.Range("C216").Value2 = "Step"
$ws.Range("D216").Value2 = "步驟"
$ws.Range("E216").Value2 = "Schritt"
$ws.Range("F216").Value2 = "步骤"

$ws.Range("B217").Value2 = "seq_sequence_setup"
$ws.Range("C217").Value2 = "Sequence Setup"
$ws.Range("D217").Value2 = "流程準備程序"
$ws.Range("E217").Value2 = "Sequenz-Setup"
$ws.Range("F217").Value2 = "流程准备程序"

$ws.Range("B218").Value2 = "seq_sequence_teardown"
$ws.Range("C218").Value2 = "Sequence Teardown"
$ws.Range("D218").Value2 = "流程結束程序"
$ws.Range("E218").Value2 = "Sequenz Teardown"
$ws.Range("F218").Value2 = "流程结束程序"

$ws.Range("B219").Value2 = "seq_enable"
$ws.Range("C219").Value2 = "Enable"
$ws.Range("D219").Value2 = "啟用"
$ws.Range("E219").Value2 = "Aktivieren"
$ws.Range("F219").Value2 = "启用"

$ws.Range("B220").Value2 = "seq_delete"
$ws.Range("C220").Value2 = "Delete"
$ws.Range("D220").Value2 = "刪除"
$ws.Range("E220").Value2 = "Löschen"
$ws.Range("F220").Value2 = "删除"

$ws.Range("B221").Value2 = "open_teach_pos_pdf"
$ws.Range("C221").Value2 = "Click me to open guide of teaching position"
$ws.Range("D221").Value2 = "點擊我打開教學位置指南"
$ws.Range("E221").Value2 = "Klicken Sie auf mich, um den Leitfaden für die Lehrposition zu öffnen"
$ws.Range("F221").Value2 = "点击我打开教学位置指南"

$ws.Range("B222").Value2 = "side_posteach"
$ws.Range("C222").Value2 = "MearPos Guide"
$ws.Range("D222").Value2 = "量測點教學"
$ws.Range("E222").Value2 = "MearPos Führer"
$ws.Range("F222").Value2 = "量测点教学"

Write-Host "appended rows complete"


